$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the filename casing in column E (Source File) for all existing data rows (2-320):
#    "Special_surgery" -> "Special_Surgery"
$ws.Range("E2:E320").Replace("Y4_B2526_General_&_Special_surgery_1_B1_reference_data.xlsx", "Y4_B2526_General_&_Special_Surgery_1_B1_reference_data.xlsx")

# 2. Append a new row 321 with the new student record.
#    Set values first (use a leading apostrophe on the numeric-looking Student ID
#    so it is stored as text, matching the existing column A formatting), then copy
#    the row-striping format from row 319 (style "4") onto the new row afterwards so
#    the format-only paste does not disturb the values we already entered.
$ws.Cells.Item(321, 1).Value = "'223007"
$ws.Cells.Item(321, 2).Value = "خالد احمد محمد الكردى"
$ws.Cells.Item(321, 3).Value = "Year 4"
$ws.Cells.Item(321, 4).Value = "B1F2"
$ws.Cells.Item(321, 5).Value = "Y4_B2526_General_&_Special_Surgery_1_B1_reference_data.xlsx"

$ws.Range("A319:E319").Copy()
$ws.Range("A321:E321").PasteSpecial(-4122)
